$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 93.967282776995887
$ws.Range("C2").Value = 93.59849480056269
$ws.Range("D2").Value = 94.821777757501053
$ws.Range("E2").Value = 95.112584176485967

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 94.067146827857201
$ws.Range("C3").Value = 97.738090329780576
$ws.Range("D3").Value = 96.43580666409791
$ws.Range("E3").Value = 95.275056922025485

# Update the selection to match the new, smaller selected range
$ws.Range("B1:E3").Select()
